# 0609: risk matrix cnstr
# Insert two new constraint columns "G" and "S" into the table (表1) right
# after column "D" and before column "N", filling the 4 data rows with the
# same values Excel wrote for the new columns ("0.1" for G, "inf" for S -
# matching the existing "N" constraint's value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Physically insert 2 blank columns at L:M (pushes N..P -> P..R), carrying
# the existing cell formatting/styles along with them.
$ws.Columns("L:M").Insert()

# Grow the table ("表1") so it spans the two new columns.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:R5"))

# Re-assert the header text for every column from the insertion point
# onward (L..R) - the two new ones plus the ones that shifted right - so
# the table's column list picks up the correct names in the correct slots.
$ws.Range("L1").Value = "G"
$ws.Range("M1").Value = "S"
$ws.Range("N1").Value = "N"
$ws.Range("O1").Value = "wei_tole"
$ws.Range("P1").Value = "begin_date"
$ws.Range("Q1").Value = "end_date"
$ws.Range("R1").Value = "opt_verbose"

# Fill the new columns' data cells (rows 2-5).
$ws.Range("L2:L5").Value = "0.1"
$ws.Range("M2:M5").Value = "inf"

# Match the author's final selection.
$ws.Range("L2").Select()
